$wb = $excel.ActiveWorkbook

$oldText = "February 03 2026 17.29.55 EST"
$newText = "February 03 2026 18.05.36 EST"

foreach ($ws in $wb.Worksheets) {
    foreach ($cell in $ws.UsedRange.Cells) {
        $v = $cell.Value2
        if ($v -ne $null -and $v -is [string] -and $v.Contains($oldText)) {
            $cell.Value2 = $v.Replace($oldText, $newText)
        }
    }
}
